$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking Price text values to remain as Text (not auto-converted to numbers)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.385.83"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").Value = "2.918.58"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "350.17"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").Value = "105.93"
$ws.Range("E6").Value = "  -5.47%  "

$ws.Range("E7").Value = "  -0.91%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -3.28%  "

$ws.Range("D10").Value = "37.67"
$ws.Range("E10").Value = "  -4.40%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").Value = "0.0846"
$ws.Range("E12").Value = "  -3.61%  "

$ws.Range("D13").Value = "18.87"
$ws.Range("E13").Value = "  -6.07%  "

$ws.Range("D14").Value = "3.375.42"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "7.57"
$ws.Range("E15").Value = "  -2.46%  "

$ws.Range("D16").Value = "2.922.47"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "0.957"
$ws.Range("E17").Value = "  -2.44%  "

$ws.Range("D18").Value = "51.371.08"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").Value = "3.40"
$ws.Range("E19").Value = "  +3.29%  "

$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("D21").Value = "13.37"
$ws.Range("E21").Value = "  -5.81%  "

$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").Value = "68.78"
$ws.Range("E23").Value = "  -3.35%  "

$ws.Range("D24").Value = "259.61"
$ws.Range("E24").Value = "  -3.21%  "

$ws.Range("E25").Value = "  -3.46%  "

$ws.Range("E26").Value = "  -3.35%  "

$ws.Range("D27").Value = "26.32"
$ws.Range("E27").Value = "  -2.51%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  +3.42%  "

$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("D31").Value = "10.20"
$ws.Range("E31").Value = "  -3.80%  "

$ws.Range("D32").Value = "6.03"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("D33").Value = "35.56"
$ws.Range("E33").Value = "  -3.95%  "

$ws.Range("E34").Value = "  -5.12%  "

$ws.Range("D35").Value = "50.33"
$ws.Range("E35").Value = "  -5.09%  "

$ws.Range("D36").Value = "0.0426"
$ws.Range("E36").Value = "  -5.95%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  -7.75%  "

$ws.Range("D39").Value = "17.58"
$ws.Range("E39").Value = "  -5.77%  "

$ws.Range("E40").Value = "  -6.19%  "

$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("D43").Value = "22.13"
$ws.Range("E43").Value = "  -5.72%  "

$ws.Range("D44").Value = "119.55"
$ws.Range("E44").Value = "  +6.69%  "

$ws.Range("E45").Value = "  -3.08%  "

$ws.Range("D46").Value = "2.091.85"
$ws.Range("E46").Value = "  -4.85%  "

$ws.Range("D47").Value = "3.29"
$ws.Range("E47").Value = "  -6.33%  "

$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  -8.92%  "

$ws.Range("E49").Value = "  -4.45%  "

$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  -5.54%  "

$ws.Range("D51").Value = "0.905"
$ws.Range("E51").Value = "  -4.95%  "

# Restore default (Normal) cell style so only the value changed, matching original formatting
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
